$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tradewheel_Leads")

# Update column widths: A 16 -> 20, B 14 -> 16
# (COM ColumnWidth has a +5/6 offset relative to the stored OOXML width)
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws.Columns.Item(2).ColumnWidth = 15.166666666666668

# Write the refreshed scrape rows (2-23)
$ws.Cells.Item(2, 1).Value = 'United Kingdom'
$ws.Cells.Item(2, 2).Value = '31 minutes ago'
$ws.Cells.Item(2, 3).Value = 'Prices For Soybean and Corn'
$ws.Cells.Item(2, 4).Value = 'https://www.tradewheel.com/buyers/prices-for-soybean-and-corn/902542/'
$ws.Cells.Item(2, 5).Value = 'Hello We are looking for Soybean and Corn in UK We need about 50MT both Soybean and Corn Kindly share your complete company profile Thank you'
$ws.Cells.Item(2, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(3, 1).Value = 'USA'
$ws.Cells.Item(3, 2).Value = '18 hours ago'
$ws.Cells.Item(3, 3).Value = 'Need prices for charger'
$ws.Cells.Item(3, 4).Value = 'https://www.tradewheel.com/buyers/need-prices-for-charger/902386/'
$ws.Cells.Item(3, 5).Value = 'I am looking for the 100 pcs of chargers that have 3 different ends on one cord. Type C, Iphone , things like that. I am a resaler trying to get a business up and running so looking for just about anything'
$ws.Cells.Item(3, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(4, 1).Value = 'Italy'
$ws.Cells.Item(4, 2).Value = '18 hours ago'
$ws.Cells.Item(4, 3).Value = 'Need to Purchase Makeup brushes, Wholesale in bulk'
$ws.Cells.Item(4, 4).Value = 'https://www.tradewheel.com/buyers/need-to-purchase-makeup-brushes-wholesale-in-bulk/902383/'
$ws.Cells.Item(4, 5).Value = 'Hi I am looking to purchase high-quality vegan makeup brushes, eyeshadow palettes, and beauty sponges. I am also interested in exploring custom-branded options and reviewing your available models to select the most suitable products for my brand....'
$ws.Cells.Item(4, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(5, 1).Value = 'Australia'
$ws.Cells.Item(5, 2).Value = '19 hours ago'
$ws.Cells.Item(5, 3).Value = 'Importing bottle – Need quotes'
$ws.Cells.Item(5, 4).Value = 'https://www.tradewheel.com/buyers/importing-bottle-need-quotes/902374/'
$ws.Cells.Item(5, 5).Value = 'Good day, We would like to purchase water bottles for school and offices use. Hoping to connect soon and explore cooperation opportunities.'
$ws.Cells.Item(5, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(6, 1).Value = 'Singapore'
$ws.Cells.Item(6, 2).Value = '19 hours ago'
$ws.Cells.Item(6, 3).Value = 'Procuring Float Glass for Building and Interior Projects'
$ws.Cells.Item(6, 4).Value = 'https://www.tradewheel.com/buyers/procuring-float-glass-for-building-and-interior-projects/902372/'
$ws.Cells.Item(6, 5).Value = 'I want to procure low iron float glass. Thickness: 10 mm Size: 1500 mm x 600 mm Quantity: 4 sheets Delivery to: Singapore Please provide FOB/CIF price and lead time'
$ws.Cells.Item(6, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(7, 1).Value = 'China'
$ws.Cells.Item(7, 2).Value = '19 hours ago'
$ws.Cells.Item(7, 3).Value = 'Buying sugar for export.'
$ws.Cells.Item(7, 4).Value = 'https://www.tradewheel.com/buyers/buying-sugar-for-export/902371/'
$ws.Cells.Item(7, 5).Value = 'Greetings, We require bulk quantities of sugar at wholesale price. Kindly share your catalog, minimum order quantity along with delivery timeline.'
$ws.Cells.Item(7, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(8, 1).Value = 'China'
$ws.Cells.Item(8, 2).Value = '19 hours ago'
$ws.Cells.Item(8, 3).Value = 'Sourcing Hair Treatment Products for Salon and Retail Use'
$ws.Cells.Item(8, 4).Value = 'https://www.tradewheel.com/buyers/sourcing-hair-treatment-products-for-salon-and-retail/902370/'
$ws.Cells.Item(8, 5).Value = 'We are purchasing hair treatment products for salons. Types: oils, serums, masks Packaging: bottles, tubes MOQ: 500 units per SKU'
$ws.Cells.Item(8, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(9, 1).Value = 'China'
$ws.Cells.Item(9, 2).Value = '19 hours ago'
$ws.Cells.Item(9, 3).Value = 'Buying Dump Truck for Heavy-Duty Transport and Logistics'
$ws.Cells.Item(9, 4).Value = 'https://www.tradewheel.com/buyers/buying-dump-truck-for-heavy-duty-transport-and-logistics/902369/'
$ws.Cells.Item(9, 5).Value = 'We need 2 units of dump trucks, heavy-duty, suitable for construction sites. Please send full specifications and quotation for review.'
$ws.Cells.Item(9, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(10, 1).Value = 'USA'
$ws.Cells.Item(10, 2).Value = '19 hours ago'
$ws.Cells.Item(10, 3).Value = 'Procuring Car Radio for Automotive Aftermarket Supply'
$ws.Cells.Item(10, 4).Value = 'https://www.tradewheel.com/buyers/procuring-car-radio-for-automotive-aftermarket-supply/902368/'
$ws.Cells.Item(10, 5).Value = 'I’m looking for car radios, 50 pcs, Bluetooth-enabled, FM/AM, suitable for standard dashboards. Please respond with your quotation and lead time.'
$ws.Cells.Item(10, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(11, 1).Value = 'United Kingdom'
$ws.Cells.Item(11, 2).Value = '19 hours ago'
$ws.Cells.Item(11, 3).Value = 'Need Korean Skincare for Cosmetic Retail Distribution'
$ws.Cells.Item(11, 4).Value = 'https://www.tradewheel.com/buyers/need-korean-skincare-for-cosmetic-retail-distribution/902367/'
$ws.Cells.Item(11, 5).Value = 'I’m searching Korean skincare products, 500 units, hydrating and anti-aging formulations, ready for retail packaging. Awaiting your commercial offer for evaluation.'
$ws.Cells.Item(11, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(12, 1).Value = 'Zambia'
$ws.Cells.Item(12, 2).Value = '19 hours ago'
$ws.Cells.Item(12, 3).Value = 'Buying Electric Heaters for Seasonal Stock and Retail Supply'
$ws.Cells.Item(12, 4).Value = 'https://www.tradewheel.com/buyers/buying-electric-heaters-for-seasonal-stock-and-retail/902366/'
$ws.Cells.Item(12, 5).Value = 'I want to source electric heaters in bulk. Features: overheat protection, thermostat Packaging: retail box or bulk Delivery: FOB/CIF'
$ws.Cells.Item(12, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(13, 1).Value = 'South Africa'
$ws.Cells.Item(13, 2).Value = '19 hours ago'
$ws.Cells.Item(13, 3).Value = 'Importing Construction Machinery Parts for Equipment Maintenance'
$ws.Cells.Item(13, 4).Value = 'https://www.tradewheel.com/buyers/importing-construction-machinery-parts-for-equipment-maintenance/902364/'
$ws.Cells.Item(13, 5).Value = 'We need 150 units of construction machinery parts, precision-made and durable, suitable for repair and assembly. Awaiting your commercial offer for evaluation.'
$ws.Cells.Item(13, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(14, 1).Value = 'Dominican Republic'
$ws.Cells.Item(14, 2).Value = '19 hours ago'
$ws.Cells.Item(14, 3).Value = 'Procuring Dining Table for Home and Office Furniture Supply'
$ws.Cells.Item(14, 4).Value = 'https://www.tradewheel.com/buyers/procuring-dining-table-for-home-and-office-furniture/902363/'
$ws.Cells.Item(14, 5).Value = 'We require dining tables for wholesale furniture supply. Material: solid wood, engineered wood, glass top Design: modern, classic Size: 120–200 cm length MOQ: 50 pcs'
$ws.Cells.Item(14, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(15, 1).Value = 'France'
$ws.Cells.Item(15, 2).Value = '19 hours ago'
$ws.Cells.Item(15, 3).Value = 'Sourcing CPUs for Computer Assembly and IT Upgrades'
$ws.Cells.Item(15, 4).Value = 'https://www.tradewheel.com/buyers/sourcing-cpus-for-computer-assembly-and-it-upgrades/902362/'
$ws.Cells.Item(15, 5).Value = 'We are purchasing CPUs for computer assembly. Type: Intel i5/i7, AMD Ryzen 5/7 Clock speed: 2.5–4.5 GHz Warranty: 1 year MOQ: 50 pcs'
$ws.Cells.Item(15, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(16, 1).Value = 'Albania'
$ws.Cells.Item(16, 2).Value = '19 hours ago'
$ws.Cells.Item(16, 3).Value = 'Looking to partner with baby food supplier'
$ws.Cells.Item(16, 4).Value = 'https://www.tradewheel.com/buyers/looking-to-partner-with-baby-food-supplier/902361/'
$ws.Cells.Item(16, 5).Value = 'Good day, We are currently sourcing baby food. The products must be made from high-quality ingredients with no harmful additives.'
$ws.Cells.Item(16, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(17, 1).Value = 'Hong Kong'
$ws.Cells.Item(17, 2).Value = '19 hours ago'
$ws.Cells.Item(17, 3).Value = 'Buying Steel Round Bars for Construction and Fabrication Projects'
$ws.Cells.Item(17, 4).Value = 'https://www.tradewheel.com/buyers/buying-steel-round-bars-for-construction-and-fabrication/902360/'
$ws.Cells.Item(17, 5).Value = 'I am looking for suppliers of steel round bars. Surface: hot-rolled, cold-rolled Grade: ASTM A36, AISI 304 Packaging: bundled or palletized'
$ws.Cells.Item(17, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(18, 1).Value = 'Germany'
$ws.Cells.Item(18, 2).Value = '19 hours ago'
$ws.Cells.Item(18, 3).Value = 'Honey needed for retain chain supply'
$ws.Cells.Item(18, 4).Value = 'https://www.tradewheel.com/buyers/honey-needed-for-retain-chain-supply/902359/'
$ws.Cells.Item(18, 5).Value = 'Good day, We are interested in buying honey. Could you provide us with availability along with pricing? Kind Regards.'
$ws.Cells.Item(18, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(19, 1).Value = 'Hong Kong'
$ws.Cells.Item(19, 2).Value = '19 hours ago'
$ws.Cells.Item(19, 3).Value = 'Required dried seafood in bulk quantities'
$ws.Cells.Item(19, 4).Value = 'https://www.tradewheel.com/buyers/required-dried-seafood-in-bulk-quantities/902357/'
$ws.Cells.Item(19, 5).Value = 'Hello, We want to purchase dried seafood suitable for retail sale and food processing. You can reach me through email for quick communication.'
$ws.Cells.Item(19, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(20, 1).Value = 'Hong Kong'
$ws.Cells.Item(20, 2).Value = '19 hours ago'
$ws.Cells.Item(20, 3).Value = 'Bulk purchase inquiry for plastic cups'
$ws.Cells.Item(20, 4).Value = 'https://www.tradewheel.com/buyers/bulk-purchase-inquiry-for-plastic-cups/902356/'
$ws.Cells.Item(20, 5).Value = 'Hello, We are currently seeking plastic cups. We prefer suppliers offering custom printing and eco-friendly options. Kind Regards.'
$ws.Cells.Item(20, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(21, 1).Value = 'Pakistan'
$ws.Cells.Item(21, 2).Value = '19 hours ago'
$ws.Cells.Item(21, 3).Value = 'Need cookware sets for distribution'
$ws.Cells.Item(21, 4).Value = 'https://www.tradewheel.com/buyers/need-cookware-sets-for-distribution/902355/'
$ws.Cells.Item(21, 5).Value = 'Greetings, We are interested in buying cookware sets. The sets must include durable materials, even heat distribution & stainless-steel surfaces.'
$ws.Cells.Item(21, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(22, 1).Value = 'India'
$ws.Cells.Item(22, 2).Value = '1 day ago'
$ws.Cells.Item(22, 3).Value = 'Importing Car Refrigerators – Need Quotes'
$ws.Cells.Item(22, 4).Value = 'https://www.tradewheel.com/buyers/importing-car-refrigerators-need-quotes/902134/'
$ws.Cells.Item(22, 5).Value = 'Hello, We would like to purchase Car Refrigerators. We require units with fast cooling, low power consumption and stable performance.'
$ws.Cells.Item(22, 6).Value = '2025-11-19 13:48:48'

$ws.Cells.Item(23, 1).Value = 'Japan'
$ws.Cells.Item(23, 2).Value = '1 day ago'
$ws.Cells.Item(23, 3).Value = 'Buying Telecentric Lens for High-Accuracy Industrial Imaging Projects'
$ws.Cells.Item(23, 4).Value = 'https://www.tradewheel.com/buyers/buying-telecentric-lens-for-high-accuracy-industrial-imaging-projects/902133/'
$ws.Cells.Item(23, 5).Value = 'We require telecentric lens options with low distortion and adjustable working distance for precision measurement tasks. Please respond with your quotation and lead time.'
$ws.Cells.Item(23, 6).Value = '2025-11-19 13:48:48'
